# Use "Then" instead of "Assert" to match BDD syntax
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the cell value from "Assert" to "Then"
$ws.Range("A12").Value = "Then"

# Update the conditional formatting rule that matched "Assert" to now match "Then"
$fcs = $ws.Range("A1:XFD1048576").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Type -eq 1 -and $fc.Formula1 -eq '="Assert"') {
        $fc.Formula1 = '="Then"'
    }
}
